# Update the LR-pair sheet with new TPM-based NATMI values and add the
# "Resolving-Mac" target cluster (rows 5,9,13 previously missing / rows
# reshuffled + 3 new rows appended for the MuSCs sending cluster).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then columns E..T (numeric statistics recomputed with the new TPM data).
$data = @(
    @("ECs", "Fgf17", "Fgfr3", "ECs", 1, 0.3333333333333333, 0.1816793333333333, 0.545038, 0.3544748187589303, 0.3544748187589304, 3, 1, 4.618552666666667, 13.855658, 0.78434648953826, 0.78434648953826, 0.839095569444889, 7.551860125004001, 0.278031079723278, 0.278031079723278),
    @("ECs", "Fgf17", "Fgfr3", "FAPs", 1, 0.3333333333333333, 0.1816793333333333, 0.545038, 0.3544748187589303, 0.3544748187589304, 2, 0.6666666666666666, 0.6792986666666666, 2.037896, 0.1153620112191035, 0.1153620112191036, 0.1234145288942222, 1.110730760048, 0.0408929280185574, 0.04089292801855743),
    @("ECs", "Fgf17", "Fgfr3", "MuSCs", 1, 0.3333333333333333, 0.1816793333333333, 0.545038, 0.3544748187589303, 0.3544748187589304, 3, 1, 0.478937, 1.436811, 0.0813355572127976, 0.08133555721279762, 0.08701295486866667, 0.783116593818, 0.02883140690166303, 0.02883140690166305),
    @("ECs", "Fgf17", "Fgfr3", "Resolving-Mac", 1, 0.3333333333333333, 0.1816793333333333, 0.545038, 0.3544748187589303, 0.3544748187589304, 1, 0.3333333333333333, 0.1116203333333333, 0.334861, 0.01895594202983873, 0.01895594202983874, 0.02027910774644445, 0.182511969718, 0.006719404115431873, 0.006719404115431877),
    @("FAPs", "Fgf17", "Fgfr3", "ECs", 1, 0.3333333333333333, 0.07524966666666667, 0.225749, 0.14681973708257, 0.14681973708257, 3, 1, 4.618552666666667, 13.855658, 0.78434648953826, 0.78434648953826, 0.3475445486491112, 3.127900937842, 0.115157545375644, 0.115157545375644),
    @("FAPs", "Fgf17", "Fgfr3", "FAPs", 1, 0.3333333333333333, 0.07524966666666667, 0.225749, 0.14681973708257, 0.14681973708257, 2, 0.6666666666666666, 0.6792986666666666, 2.037896, 0.1153620112191035, 0.1153620112191036, 0.05111699823377778, 0.460052984104, 0.01693742015650527, 0.01693742015650527),
    @("FAPs", "Fgf17", "Fgfr3", "MuSCs", 1, 0.3333333333333333, 0.07524966666666667, 0.225749, 0.14681973708257, 0.14681973708257, 3, 1, 0.478937, 1.436811, 0.0813355572127976, 0.08133555721279762, 0.03603984960433334, 0.324358646439, 0.01194166512544727, 0.01194166512544727),
    @("FAPs", "Fgf17", "Fgfr3", "Resolving-Mac", 1, 0.3333333333333333, 0.07524966666666667, 0.225749, 0.14681973708257, 0.14681973708257, 1, 0.3333333333333333, 0.1116203333333333, 0.334861, 0.01895594202983873, 0.01895594202983874, 0.008399392876555556, 0.075594535889, 0.00278310642497336, 0.002783106424973361),
    @("MuSCs", "Fgf17", "Fgfr3", "ECs", 3, 1, 0.255602, 0.766806, 0.4987054441584996, 0.4987054441584997, 3, 1, 4.618552666666667, 13.855658, 0.78434648953826, 0.78434648953826, 1.180511298705333, 10.624601688348, 0.3911578644393379, 0.391157864439338),
    @("MuSCs", "Fgf17", "Fgfr3", "FAPs", 3, 1, 0.255602, 0.766806, 0.4987054441584996, 0.4987054441584997, 2, 0.6666666666666666, 0.6792986666666666, 2.037896, 0.1153620112191035, 0.1153620112191036, 0.1736300977973333, 1.562670880176, 0.05753166304404084, 0.05753166304404086),
    @("MuSCs", "Fgf17", "Fgfr3", "MuSCs", 3, 1, 0.255602, 0.766806, 0.4987054441584996, 0.4987054441584997, 3, 1, 0.478937, 1.436811, 0.0813355572127976, 0.08133555721279762, 0.122417255074, 1.101755295666, 0.04056248518568728, 0.04056248518568729),
    @("MuSCs", "Fgf17", "Fgfr3", "Resolving-Mac", 3, 1, 0.255602, 0.766806, 0.4987054441584996, 0.4987054441584997, 1, 0.3333333333333333, 0.1116203333333333, 0.334861, 0.01895594202983873, 0.01895594202983874, 0.02853038044066667, 0.256773423966, 0.009453431489433495, 0.009453431489433499)
)

$r = 2
foreach ($row in $data) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $row[$i]
    }
    $r = $r + 1
}
